$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three footnote cells with corrected wording.
$ws.Range("A5").Value = "(Half-rate F32 Accumulate during training)"
$ws.Range("A6").Value = "(Half the tensor cores are split between 16/32)"
$ws.Range("A7").Value = "(Poor drivers. Might not impact tensor core performance for training.)"

# Move the two reference URLs down two rows (A17->A19, A18 cleared) and
# relocate the "Note" text that used to sit at A20 down to A22, with the
# exxactcorp URL taking its old spot at A20.
$evolutionUrl = $ws.Range("A17").Value2
$exxactUrl = $ws.Range("A18").Value2
$noteText = $ws.Range("A20").Value2

$ws.Range("A17").ClearContents()
$ws.Range("A18").ClearContents()

$ws.Range("A19").Value = $evolutionUrl
$ws.Range("A20").Value = $exxactUrl
$ws.Range("A22").Value = $noteText

# Update selection to match the saved view.
$ws.Range("A17").Select()
